# Update data: 11 September 2021
# Appends the 2021-08-01 (serial 44409) observations to both sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Canada": one new row (national figure)
# ---------------------------------------------------------------------
$wsCanada = $wb.Worksheets.Item("Canada")

$rowCanada = 21
$wsCanada.Cells.Item($rowCanada, 1).NumberFormat = $wsCanada.Cells.Item($rowCanada - 1, 1).NumberFormat
$wsCanada.Cells.Item($rowCanada, 1).Value = 44409
$wsCanada.Cells.Item($rowCanada, 2).NumberFormat = $wsCanada.Cells.Item($rowCanada - 1, 2).NumberFormat
$wsCanada.Cells.Item($rowCanada, 2).Value = "Canada"
$wsCanada.Cells.Item($rowCanada, 4).Value = 1440
$wsCanada.Cells.Item($rowCanada, 5).Value = 1176.5999999999999
$wsCanada.Cells.Item($rowCanada, 3).Formula = "=(D" + $rowCanada + "-E" + $rowCanada + ")/E" + $rowCanada + "*100"

# ---------------------------------------------------------------------
# Sheet "Province": ten new rows (one per province, same date)
# ---------------------------------------------------------------------
$wsProvince = $wb.Worksheets.Item("Province")

$provinceRows = @(
    @{ Row = 192; Name = "Newfoundland & Labrador"; D = 30.5;   E = 33.2 },
    @{ Row = 193; Name = "Prince Edward Island";     D = 9.4;   E = 7.5 },
    @{ Row = 194; Name = "Nova Scotia";              D = 39.5;  E = 40.4 },
    @{ Row = 195; Name = "New Brunswick";            D = 36.6;  E = 34.3 },
    @{ Row = 196; Name = "Quebec";                   D = 265.1; E = 223.1 },
    @{ Row = 197; Name = "Ontario";                  D = 607.5; E = 444.4 },
    @{ Row = 198; Name = "Manitoba";                 D = 39.8;  E = 37.9 },
    @{ Row = 199; Name = "Saskatchewan";             D = 42.3;  E = 31.9 },
    @{ Row = 200; Name = "Alberta";                  D = 193.4; E = 182.4 },
    @{ Row = 201; Name = "British Columbia";         D = 175.9; E = 141.5 }
)

foreach ($entry in $provinceRows) {
    $r = $entry.Row

    $wsProvince.Cells.Item($r, 1).NumberFormat = $wsProvince.Cells.Item($r - 1, 1).NumberFormat
    $wsProvince.Cells.Item($r, 1).Value = 44409

    if ($r -eq 192) {
        # First province of the new date group keeps the same style as the
        # analogous first-of-group cell above it (e.g. B182).
        $wsProvince.Cells.Item($r, 2).NumberFormat = $wsProvince.Cells.Item(182, 2).NumberFormat
    }
    $wsProvince.Cells.Item($r, 2).Value = $entry.Name

    $wsProvince.Cells.Item($r, 4).Value = $entry.D
    $wsProvince.Cells.Item($r, 5).Value = $entry.E
    $wsProvince.Cells.Item($r, 3).Formula = "=(D" + $r + "-E" + $r + ")/E" + $r + "*100"
}

# ---------------------------------------------------------------------
# Refresh selection / scroll position to mirror the appended rows.
# ---------------------------------------------------------------------
$wsCanada.Activate()
$wsCanada.Range("A21").Select()

$wsProvince.Activate()
$wsProvince.Range("D202").Select()

$wb.Save()
